# Update Name of Algo
# Applies updated numeric results to the KNN imputation result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value  = 13.113
$ws.Range("C3").Value  = -12.355
$ws.Range("A4").Value  = -21.468
$ws.Range("B4").Value  = 6.962999999999999
$ws.Range("C4").Value  = -12.585
$ws.Range("B5").Value  = 6.174
$ws.Range("E5").Value  = 12.769
$ws.Range("A6").Value  = -21.336
$ws.Range("B6").Value  = 6.252000000000001
$ws.Range("A7").Value  = -21.179
$ws.Range("A8").Value  = -21.398
$ws.Range("B8").Value  = 6.153
$ws.Range("C9").Value  = -11.775
$ws.Range("C11").Value = -12.642
$ws.Range("C14").Value = -11.607
$ws.Range("A16").Value = -21.212
$ws.Range("B16").Value = 5.896
$ws.Range("C18").Value = -12.621
$ws.Range("A20").Value = -22.104
$ws.Range("E20").Value = 13.33
$ws.Range("A21").Value = -21.14
$ws.Range("B22").Value = 6.386000000000001
$ws.Range("C25").Value = -12.577
